$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh -- GitHub Actions scheduled update.
# Columns: D = Price (text, locale-formatted), E = Volume(1h) change (text, padded %).

$ws.Range('D2').Value = '57.328.10'
$ws.Range('E2').Value = '  -2.85%  '
$ws.Range('D3').Value = '2.423.33'
$ws.Range('E3').Value = '  -3.52%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '511.65'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '131.02'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.52%  '
$ws.Range('E7').Value = '  +0.14%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.550'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.14%  '
$ws.Range('D9').Value = '2.424.48'
$ws.Range('E9').Value = '  -3.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0956'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.73%  '
$ws.Range('E11').Value = '  -1.50%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.21'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.84%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.330'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.47%  '
$ws.Range('D14').Value = '2.856.92'
$ws.Range('E14').Value = '  -3.37%  '
$ws.Range('D15').Value = '57.269.12'
$ws.Range('E15').Value = '  -2.79%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '21.51'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -5.72%  '
$ws.Range('E17').Value = '  -4.59%  '
$ws.Range('D18').Value = '2.437.64'
$ws.Range('E18').Value = '  -3.19%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.35'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -6.32%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '314.23'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.08'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.98%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.60'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.33%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.72'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.21%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.401'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.38%  '
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.159'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.34%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.18'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.08%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '169.35'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.09%  '
$ws.Range('D30').Value = '0.0₃0723'
$ws.Range('E30').Value = '  -5.41%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.19'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.66'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.06%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.14'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.25%  '
$ws.Range('E34').Value = '  +0.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.67'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.28'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.85'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.86%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '36.23'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.33%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.43'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -5.94%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.773'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.54%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.35'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.40%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '265.73'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.91'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.580'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '121.53'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.20%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0899'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0480'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0209'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -4.39%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '16.50'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.62%  '
$ws.Range('D51').Value = '1.702.12'
$ws.Range('E51').Value = '  -3.32%  '
